$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": the handoff for file
# 06fe9af5-eccb-4954-a050-6114de2789d0 was just (re)generated, so its
# "Latest Handoff Datetime" moves forward. The same handoff-datetime value
# is used for the row above it (2029cd60-ae90-4ac8-b96b-49efe63e8fa1, row 8)
# because that file's handoff transform is chained/depends on this one.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-09 06:35:21"
$wsZhCn.Range("D8").Value = "2016-03-09 06:35:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-09 06:35:31"
$wsDeDe.Range("D8").Value = "2016-03-09 06:35:31"
